# Write operation on excel sheet for patient details

$wb = $excel.ActiveWorkbook

# --- Sheet1: remove the last two username/password rows (Divya, Valli) ---
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Rows("6:7").Delete()
[void]$ws1.Range("B8").Select()

# --- Sheet2: build the patient-details header + data row ---
$ws2 = $wb.Worksheets.Item(2)

# Header row (the Email header is entered last, after the rest of the row,
# which is why its string lands at the end of the header block in the
# shared-string table)
$ws2.Range("A1").Value = "First Name"
$ws2.Range("B1").Value = "Middle Name"
$ws2.Range("C1").Value = "Last Name"
$ws2.Range("D1").Value = "Parent First Name"
$ws2.Range("E1").Value = "Parent Middle Name"
$ws2.Range("F1").Value = "Parent Last Name"
$ws2.Range("H1").Value = "Phone"
$ws2.Range("I1").Value = "Mail"
$ws2.Range("J1").Value = "Age"
$ws2.Range("K1").Value = "Gender"
$ws2.Range("G1").Value = "Email"

# Data row
$ws2.Range("A2").Value = "Geethika"
$ws2.Range("B2").Value = "middle name"
$ws2.Range("C2").Value = "Nannapaneni"
$ws2.Range("D2").Value = "Srinivas"
$ws2.Range("E2").Value = "Parent middle name"
$ws2.Range("F2").Value = "Nannapaneni"
$ws2.Range("G2").Value = "gnannap1@asu.edu"

# Phone/Age need to be stored as literal TEXT (not auto-converted numbers),
# so build the text through a helper formula cell and paste only the value
# across -- this keeps the digit strings as shared-string text without
# ever touching a cell's number format.
$ws2.Range("Z1").Formula = "=""6026219441"""
[void]$ws2.Range("Z1").Copy()
[void]$ws2.Range("H2").Select()
[void]$ws2.PasteSpecial("xlPasteValues")

$ws2.Range("I2").Value = "gnannap1@asu.edu"

$ws2.Range("Z1").Formula = "=""25"""
[void]$ws2.Range("Z1").Copy()
[void]$ws2.Range("J2").Select()
[void]$ws2.PasteSpecial("xlPasteValues")

$ws2.Range("K2").Value = "Female"

# Clean up the helper cell
[void]$ws2.Range("Z1").Clear()

# Make Sheet2 the active sheet/tab and select the cell below the data
[void]$ws2.Activate()
[void]$ws2.Range("H9").Select()
